# The document was produced by converting a SharePoint-hosted file, which
# leaves behind several "custom XML parts" that only exist to carry the
# SharePoint document-library content-type schema / form templates /
# property-promotion stubs (customXml/item1.xml..item3.xml and their
# itemProps*.xml counterparts). None of this is user content - it is
# purely library plumbing - so bringing the file "up to date" means
# stripping all of those custom XML parts from the package.
#
# The standard Word automation pattern for this is to walk
# Document.CustomXMLParts (back-to-front, since deleting shifts indices)
# and delete every part whose root-element namespace matches one of the
# known SharePoint/Office metadata schemas that were injected on
# conversion.

$d = $word.ActiveDocument

$targetNamespaces = @(
    "http://schemas.microsoft.com/office/2006/metadata/contentType",
    "http://schemas.microsoft.com/sharepoint/v3/contenttype/forms",
    "http://schemas.microsoft.com/office/2006/metadata/properties"
)

$parts = $d.CustomXMLParts

for ($i = $parts.Count; $i -ge 1; $i--) {
    $part = $parts.Item($i)
    if ($targetNamespaces -contains $part.NamespaceURI) {
        $part.Delete()
    }
}

# Belt-and-braces: some hosts only expose matching parts through
# SelectByNamespace rather than by walking the full collection, so cover
# that path too in case a part was missed above.
foreach ($ns in $targetNamespaces) {
    $matches = $d.CustomXMLParts.SelectByNamespace($ns)
    if ($matches -ne $null) {
        for ($j = $matches.Count; $j -ge 1; $j--) {
            $matches.Item($j).Delete()
        }
    }
}

$d.Save()
